$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '36.517.05'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").Value = '1.950.63'
$ws.Range("E3").Value = '  +0.46%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.14'
$ws.Range("E5").Value = '  -0.14%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("E6").Value = '  +0.33%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.56'
$ws.Range("E7").Value = '  +5.56%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.376'
$ws.Range("E9").Value = '  +4.01%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0787'
$ws.Range("E10").Value = '  -7.21%  '

# Row 11
$ws.Range("E11").Value = '  +0.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.24'
$ws.Range("E12").Value = '  +6.12%  '

# Row 13
$ws.Range("D13").Value = '2.235.86'
$ws.Range("E13").Value = '  +0.39%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.824'
$ws.Range("E14").Value = '  +1.53%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.54'
$ws.Range("E15").Value = '  +0.89%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.24'
$ws.Range("E16").Value = '  +1.37%  '

# Row 17
$ws.Range("D17").Value = '1.952.82'
$ws.Range("E17").Value = '  +0.93%  '

# Row 18
$ws.Range("D18").Value = '36.427.06'
$ws.Range("E18").Value = '  -0.06%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.29'
$ws.Range("E19").Value = '  -0.11%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0848'
$ws.Range("E20").Value = '  -1.79%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.89'
$ws.Range("E21").Value = '  +0.40%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.06'
$ws.Range("E22").Value = '  +1.50%  '

# Row 23
$ws.Range("E23").Value = '  +0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("E24").Value = '  +3.56%  '

# Row 25
$ws.Range("E25").Value = '  +2.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.144'
$ws.Range("E26").Value = '  +8.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.15'
$ws.Range("E27").Value = '  -0.58%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.76'
$ws.Range("E28").Value = '  -0.67%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.26'
$ws.Range("E29").Value = '  +0.42%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.31'
$ws.Range("E30").Value = '  +19.39%  '

# Row 31
$ws.Range("E31").Value = '  +0.95%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.73'
$ws.Range("E32").Value = '  +3.22%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0611'
$ws.Range("E33").Value = '  -0.58%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.43'
$ws.Range("E34").Value = '  +5.84%  '

# Row 35
$ws.Range("E35").Value = '  +9.74%  '

# Row 36
$ws.Range("E36").Value = '  -0.12%  '

# Row 37
$ws.Range("E37").Value = '  +4.56%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.76'
$ws.Range("E38").Value = '  -1.64%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.44'
$ws.Range("E39").Value = '  -12.79%  '

# Row 40
$ws.Range("B40").Value = 'Cronos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0962'
$ws.Range("E40").Value = '  -2.09%  '

# Row 41
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").Value = '  +0.61%  '

# Row 42
$ws.Range("E42").Value = '  +1.71%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0209'
$ws.Range("E43").Value = '  +0.13%  '

# Row 44
$ws.Range("D44").Value = '1.356.57'
$ws.Range("E44").Value = '  +1.23%  '

# Row 45
$ws.Range("E45").Value = '  -2.06%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.40'
$ws.Range("E46").Value = '  +2.34%  '

# Row 47
$ws.Range("E47").Value = '  -0.50%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.10'
$ws.Range("E48").Value = '  -1.29%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.83'
$ws.Range("E49").Value = '  +0.39%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.28'
$ws.Range("E50").Value = '  +4.74%  '

# Row 51
$ws.Range("D51").Value = '2.131.81'
$ws.Range("E51").Value = '  +0.63%  '
